$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1756.5714
$ws.Range("I98").Value = 2016
$ws.Range("K98").Value = 2016
$ws.Range("M98").Value = -518
$ws.Range("H107").Value = 2725.4
$ws.Range("I107").Value = 906.75
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 906.75
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 1013.25
$ws.Range("N107").Value = -13840
$ws.Range("H111").Value = 2210.1177
$ws.Range("I111").Value = 888.5454999999999
$ws.Range("K111").Value = 2665.6365
$ws.Range("M111").Value = 401.3635000000004
$ws.Range("H122").Value = 1756.5714
$ws.Range("I122").Value = 2016
$ws.Range("K122").Value = 6048
$ws.Range("M122").Value = -3598
$ws.Range("H131").Value = 4548300
$ws.Range("I131").Value = 575
$ws.Range("J131").Value = 7580117
$ws.Range("K131").Value = 1725
$ws.Range("L131").Value = 22740351
$ws.Range("M131").Value = 3315
$ws.Range("N131").Value = -22750431
$ws.Range("H138").Value = 10336.96
$ws.Range("J138").Value = 10774.067
$ws.Range("L138").Value = 32322.201
$ws.Range("N138").Value = -42602.201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19550.666
$ws.Range("I32").Value = 19550.666
$ws.Range("K32").Value = 19550.666
$ws.Range("M32").Value = -19263.666
$ws.Range("H45").Value = 3416.8948
$ws.Range("I45").Value = 2539.3333
$ws.Range("J45").Value = 4921.2856
$ws.Range("K45").Value = 2539.3333
$ws.Range("L45").Value = 4921.2856
$ws.Range("M45").Value = -2162.3333
$ws.Range("N45").Value = -5675.2856
$ws.Range("H88").Value = 3797.6155
$ws.Range("I88").Value = 2384.25
$ws.Range("J88").Value = 4425.778
$ws.Range("K88").Value = 2384.25
$ws.Range("L88").Value = 4425.778
$ws.Range("M88").Value = -1978.25
$ws.Range("N88").Value = -5237.778
$ws.Range("H91").Value = 3797.6155
$ws.Range("I91").Value = 2384.25
$ws.Range("J91").Value = 4425.778
$ws.Range("K91").Value = 2384.25
$ws.Range("L91").Value = 4425.778
$ws.Range("M91").Value = -980.25
$ws.Range("N91").Value = -7233.778
$ws.Range("H110").Value = 8745.583000000001
$ws.Range("I110").Value = 11816.667
$ws.Range("K110").Value = 11816.667
$ws.Range("M110").Value = -9771.666999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1621505.9
$ws.Range("I86").Value = 3503706.5
$ws.Range("J86").Value = 8191.143
$ws.Range("K86").Value = 3503706.5
$ws.Range("L86").Value = 8191.143
$ws.Range("M86").Value = -3502583.5
$ws.Range("N86").Value = -10437.143
$ws.Range("H89").Value = 1621505.9
$ws.Range("I89").Value = 3503706.5
$ws.Range("J89").Value = 8191.143
$ws.Range("K89").Value = 17518532.5
$ws.Range("L89").Value = 40955.715
$ws.Range("M89").Value = -17512916.5
$ws.Range("N89").Value = -52187.715
$ws.Range("H132").Value = 122395.38
$ws.Range("J132").Value = 122395.38
$ws.Range("L132").Value = 122395.38
$ws.Range("N132").Value = -132515.38
$ws.Range("H134").Value = 7145915.5
$ws.Range("I134").Value = 2839.2
$ws.Range("K134").Value = 8517.599999999999
$ws.Range("M134").Value = -5982.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34487492
$ws.Range("I31").Value = 58826870
$ws.Range("K31").Value = 58826870
$ws.Range("M31").Value = -58826575
$ws.Range("H34").Value = 34487492
$ws.Range("I34").Value = 58826870
$ws.Range("K34").Value = 58826870
$ws.Range("M34").Value = -58826668
$ws.Range("H58").Value = 4268.2666
$ws.Range("I58").Value = 2482.8572
$ws.Range("K58").Value = 2482.8572
$ws.Range("M58").Value = -2279.8572
$ws.Range("H103").Value = 43394.57
$ws.Range("J103").Value = 54594
$ws.Range("L103").Value = 54594
$ws.Range("N103").Value = -56938
$ws.Range("H107").Value = 1027.6666
$ws.Range("I107").Value = 430.41177
$ws.Range("K107").Value = 430.41177
$ws.Range("M107").Value = 1489.58823
$ws.Range("H134").Value = 1412.0834
$ws.Range("I134").Value = 1094.5
$ws.Range("K134").Value = 3283.5
$ws.Range("M134").Value = -748.5
$ws.Range("H136").Value = 4268.2666
$ws.Range("I136").Value = 2482.8572
$ws.Range("K136").Value = 7448.571599999999
$ws.Range("M136").Value = -4898.571599999999
$ws.Range("H141").Value = 592211.9
$ws.Range("J141").Value = 639896.2
$ws.Range("L141").Value = 639896.2
$ws.Range("N141").Value = -650256.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 988.14703
$ws.Range("I5").Value = 557.1667
$ws.Range("J5").Value = 1223.2273
$ws.Range("K5").Value = 1671.5001
$ws.Range("L5").Value = 3669.6819
$ws.Range("M5").Value = -1559.5001
$ws.Range("N5").Value = -3893.6819
$ws.Range("H113").Value = 2545.8667
$ws.Range("I113").Value = 2249.7778
$ws.Range("J113").Value = 2990
$ws.Range("K113").Value = 6749.3334
$ws.Range("L113").Value = 8970
$ws.Range("M113").Value = -4579.3334
$ws.Range("N113").Value = -13310
$ws.Range("H134").Value = 11108.4
$ws.Range("I134").Value = 4848
$ws.Range("J134").Value = 20499
$ws.Range("K134").Value = 14544
$ws.Range("L134").Value = 61497
$ws.Range("M134").Value = -9474
$ws.Range("N134").Value = -71637
$ws.Range("H135").Value = 988.14703
$ws.Range("I135").Value = 557.1667
$ws.Range("J135").Value = 1223.2273
$ws.Range("K135").Value = 5014.5003
$ws.Range("L135").Value = 11009.0457
$ws.Range("M135").Value = -2479.5003
$ws.Range("N135").Value = -16079.0457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 564.625
$ws.Range("I107").Value = 373.2
$ws.Range("K107").Value = 373.2
$ws.Range("M107").Value = 1546.8
$ws.Range("H134").Value = 51142.43
$ws.Range("J134").Value = 51142.43
$ws.Range("L134").Value = 153427.29
$ws.Range("N134").Value = -158497.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2552.6667
$ws.Range("J22").Value = 3875.6
$ws.Range("L22").Value = 3875.6
$ws.Range("N22").Value = -4465.6
$ws.Range("H27").Value = 2552.6667
$ws.Range("J27").Value = 3875.6
$ws.Range("L27").Value = 3875.6
$ws.Range("N27").Value = -4089.6
$ws.Range("H55").Value = 1159.6136
$ws.Range("I55").Value = 823.069
$ws.Range("J55").Value = 1810.2667
$ws.Range("K55").Value = 823.069
$ws.Range("L55").Value = 1810.2667
$ws.Range("M55").Value = -650.069
$ws.Range("N55").Value = -2156.2667
$ws.Range("H68").Value = 3603540.8
$ws.Range("I68").Value = 4421145.5
$ws.Range("K68").Value = 4421145.5
$ws.Range("M68").Value = -4420396.5
$ws.Range("H71").Value = 3603540.8
$ws.Range("I71").Value = 4421145.5
$ws.Range("K71").Value = 22105727.5
$ws.Range("M71").Value = -22101983.5
$ws.Range("H135").Value = 65000
$ws.Range("J135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140
$ws.Range("H136").Value = 10065.385
$ws.Range("I136").Value = 15415.333
$ws.Range("J136").Value = 5479.7144
$ws.Range("K136").Value = 46245.999
$ws.Range("L136").Value = 16439.1432
$ws.Range("M136").Value = -43695.999
$ws.Range("N136").Value = -21539.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 627953.6
$ws.Range("I132").Value = 3150.5334
$ws.Range("K132").Value = 9451.600199999999
$ws.Range("M132").Value = -6921.600199999999
$ws.Range("H136").Value = 424149.38
$ws.Range("I136").Value = 7708.409
$ws.Range("K136").Value = 23125.227
$ws.Range("M136").Value = -20575.227
